# M07 Frozen Token Embeddings — refresh the per-epoch accuracy table on Sheet1.
# Column A holds the epoch index (0-based); for epochs >= 100 the source
# notebook's export accidentally dumped the repr() of a DisplayOutputs
# object instead of the integer, which we reproduce verbatim. Column B
# holds the validation accuracy for that epoch. The refreshed run has two
# fewer epochs than before, so rows 119-120 are dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(
    0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28,
    29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54,
    55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80,
    81, 82, 83, 84, 85, 86, 87, 88, 89, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99,
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>",
    "<__main__.DisplayOutputs object at 0x7fdfec498310>"
)

$colB = @(
    0.9375, 0.875, 0.84375, 0.734375, 0.765625, 0.734375, 0.671875, 0.65625, 0.625, 0.59375, 0.609375,
    0.53125, 0.65625, 0.625, 0.578125, 0.546875, 0.5, 0.53125, 0.578125, 0.625, 0.453125, 0.546875, 0.5,
    0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375,
    0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375,
    0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.484375, 0.46875, 0.46875, 0.46875, 0.46875,
    0.46875, 0.46875, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.453125, 0.46875, 0.46875,
    0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875,
    0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875,
    0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875, 0.46875,
    0.46875, 0.46875, 0.46875, 0.53125, 0.609375, 0.46875, 0.4375, 0.40625, 0.5625, 0.5625, 0.484375,
    0.53125, 0.546875, 0.59375, 0.578125, 0.453125, 0.484375, 0.5, 0.4262295081967213
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# The new run only has 117 epochs (rows 2-118); remove the old trailing rows.
$ws.Rows("119:120").Delete()

# Reflect the author's cursor position moving from I16 to I15.
$ws.Range("A1:XFD1048576").Select()
$ws.Range("I15").Activate()
